## Generate Report for Handback
##
## This script mirrors the "localization-status.xlsx" handback report refresh:
##  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
##    (shown on the Overview sheet as well as the zh-cn / de-de detail sheets).
##  - The per-language detail sheets (zh-cn, de-de) gain "Latest Target File" (F)
##    and "Latest Handback File" (G) entries (file name + hyperlink), mirroring
##    the existing "Source File Name" (A) / "Latest Handoff File" (D) links.
##  - "Latest Handback DateTime" (H) is stamped with the handback timestamp for
##    each language.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Update Status column -------------------------------------------------
# Overview sheet: column B = zh-cn status, column C = de-de status.
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# Detail sheets: column C = Status.
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Helper to add a "Latest Target File" / "Latest Handback File" style --
# file-name + hyperlink cell, matching the look of the existing hyperlinked
# cells (A/B/D columns) on the detail sheets.
function Add-HandbackFileLink {
    param($ws, $cellAddr, $text, $url)

    $rng = $ws.Range($cellAddr)
    $rng.Value = $text
    $rng.Font.Underline = $true
    $rng.Font.Color = 15570276
    $ws.Hyperlinks.Add($rng, $url, "", "", $text) | Out-Null
}

# --- zh-cn sheet: columns F (Latest Target File) / G (Latest Handback File)
Add-HandbackFileLink $wsZhCn "F2" "3d898684-bb87-4b1c-9483-382e80c605f1.md" "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/3d898684-bb87-4b1c-9483-382e80c605f1.md"
Add-HandbackFileLink $wsZhCn "G2" "3d898684-bb87-4b1c-9483-382e80c605f1.6d7659a2901e333762a357162aed2fcfd0b2a5e8.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e3a1c113d8d2aef05bd6809a2d0157ed7063af5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/3d898684-bb87-4b1c-9483-382e80c605f1.6d7659a2901e333762a357162aed2fcfd0b2a5e8.zh-cn.xlf"
Add-HandbackFileLink $wsZhCn "F3" "56e3ed2d-1133-457c-be66-a1c4ccf2e580.md" "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/56e3ed2d-1133-457c-be66-a1c4ccf2e580.md"
Add-HandbackFileLink $wsZhCn "G3" "56e3ed2d-1133-457c-be66-a1c4ccf2e580.88f1c291211ef5ec4d866582e427fdf0faa75d1c.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e3a1c113d8d2aef05bd6809a2d0157ed7063af5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/56e3ed2d-1133-457c-be66-a1c4ccf2e580.88f1c291211ef5ec4d866582e427fdf0faa75d1c.zh-cn.xlf"

# --- de-de sheet: columns F (Latest Target File) / G (Latest Handback File)
Add-HandbackFileLink $wsDeDe "F2" "3d898684-bb87-4b1c-9483-382e80c605f1.md" "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/3d898684-bb87-4b1c-9483-382e80c605f1.md"
Add-HandbackFileLink $wsDeDe "G2" "3d898684-bb87-4b1c-9483-382e80c605f1.6d7659a2901e333762a357162aed2fcfd0b2a5e8.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/65d07f5e9a779efb85ea7850fe0f289c136c2e69/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/3d898684-bb87-4b1c-9483-382e80c605f1.6d7659a2901e333762a357162aed2fcfd0b2a5e8.de-de.xlf"
Add-HandbackFileLink $wsDeDe "F3" "56e3ed2d-1133-457c-be66-a1c4ccf2e580.md" "https://github.com/OpenLocalizationTest/oltest/blob/d4ed5db2b70961318fb19e8f4519265ba553db52/e2e/56e3ed2d-1133-457c-be66-a1c4ccf2e580.md"
Add-HandbackFileLink $wsDeDe "G3" "56e3ed2d-1133-457c-be66-a1c4ccf2e580.88f1c291211ef5ec4d866582e427fdf0faa75d1c.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/65d07f5e9a779efb85ea7850fe0f289c136c2e69/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/56e3ed2d-1133-457c-be66-a1c4ccf2e580.88f1c291211ef5ec4d866582e427fdf0faa75d1c.de-de.xlf"

# --- Latest Handback DateTime (column H) -----------------------------------
$wsZhCn.Range("H2").Value = "2016-03-17 03:24:08"
$wsZhCn.Range("H3").Value = "2016-03-17 03:24:08"
$wsDeDe.Range("H2").Value = "2016-03-17 03:24:22"
$wsDeDe.Range("H3").Value = "2016-03-17 03:24:22"

Write-Host "Generated handback report."
